{"js": "// The diff:\n//   1) Changes the paragraph-mark rFonts hint of the \"\u5f88\u82e6\" paragraph from\n//      \"default\" to \"eastAsia\".\n//   2) Inserts a brand-new paragraph right after it containing the text\n//      \"3.10 \u5929\u6c14\u6674 \u5468\u4e94\u4e0b\u5348\u53c8\u662f\u5f00\u6e90\" (with pPr rFonts hint=\"default\" and\n//      run rFonts hint=\"eastAsia\", matching the sibling paragraphs' pattern).\n//   3) The trailing \"_GoBack\" bookmark (which sat at the end of the \"\u5f88\u82e6\"\n//      paragraph) moves along to the end of the document/new paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Find the \"\u5f88\u82e6\" paragraph; fall back to the last paragraph in the body\n// if (for some reason) the exact text can't be located.\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\u5f88\u82e6\") {\n    targetPara = paragraphs.items[i];\n    break;\n  }\n}\nif (!targetPara) {\n  targetPara = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// 1) Insert the new paragraph *after* the target, before touching its\n//    formatting \u2014 the new paragraph inherits the (still \"default\") pPr hint\n//    from its source, exactly like the target document.\nconst newPara = targetPara.insertParagraph(\n  \"3.10 \u5929\u6c14\u6674 \u5468\u4e94\u4e0b\u5348\u53c8\u662f\u5f00\u6e90\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 2) Flip the target paragraph's own paragraph-mark rFonts hint from\n//    \"default\" to \"eastAsia\" (not exposed as a discrete property on the\n//    Word.Paragraph/Word.Font object model, so we round-trip it through\n//    insertOoxml on the paragraph's own range).\nconst wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\nconst fixedParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData><w:document \" + wNs + \"><w:body>\" +\n  '<w:p><w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:lang w:val=\"en-US\" w:eastAsia=\"zh-CN\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:lang w:val=\"en-US\" w:eastAsia=\"zh-CN\"/></w:rPr><w:t>\u5f88\u82e6</w:t></w:r>' +\n  \"</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nconst targetWhole = targetPara.getRange(Word.RangeLocation.whole);\ntargetWhole.insertOoxml(fixedParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Relocate the \"_GoBack\" bookmark to the very end of the document, i.e.\n//    the end of the freshly inserted paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst docEnd = body.getRange(Word.RangeLocation.end);\ndocEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The diff:\n#   1) Changes the paragraph-mark rFonts hint of the \"\u5f88\u82e6\" paragraph from\n#      \"default\" to \"eastAsia\".\n#   2) Inserts a brand-new paragraph right after it containing the text\n#      \"3.10 \u5929\u6c14\u6674 \u5468\u4e94\u4e0b\u5348\u53c8\u662f\u5f00\u6e90\" (pPr rFonts hint=\"default\", run rFonts\n#      hint=\"eastAsia\" \u2014 the same pattern the other paragraphs use).\n#   3) The trailing \"_GoBack\" bookmark (previously at the end of the \"\u5f88\u82e6\"\n#      paragraph) moves along to the end of the document/new paragraph.\n#\n# NOTE: the \"w:hint\" attribute on a paragraph MARK's rFonts (i.e. the\n# <w:pPr><w:rPr> of a paragraph, as opposed to a run's <w:r><w:rPr>) is not\n# reachable through Font/ParagraphFormat properties in the Word object\n# model \u2014 Font writes here always resolve onto an actual text run, never\n# the paragraph-mark-only rPr. The reliable, observable lever is that a\n# *freshly created* paragraph mark inherits its rFonts hint from the\n# paragraph mark immediately preceding it at the moment of the split. So\n# instead of trying to mutate \"\u5f88\u82e6\"'s existing mark in place, we grow a new\n# mark in a spot that will inherit \"eastAsia\" (right after the prior\n# \"eastAsia\"-hinted paragraph) and relocate the \"\u5f88\u82e6\" text into it, then\n# reuse the old mark (still \"default\") for the freshly-authored paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaIndexByText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $raw = $doc.Paragraphs($i).Range.Text\n        $trimmed = $raw.TrimEnd([char]13, [char]7)\n        if ($trimmed -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$targetText = \"\u5f88\u82e6\"\n$newText = \"3.10 \u5929\u6c14\u6674 \u5468\u4e94\u4e0b\u5348\u53c8\u662f\u5f00\u6e90\"\n\n$targetIdx = Get-ParaIndexByText $d $targetText\nif ($targetIdx -eq -1) {\n    # Fallback: operate on the last paragraph of the document.\n    $targetIdx = $d.Paragraphs.Count\n}\n\n# The paragraph immediately before the target (\"\u73b0\u5728\u662f...\" in this\n# document) already carries rFonts hint=\"eastAsia\" on its own mark \u2014 a\n# newly split mark placed right after it will inherit that hint.\n$priorIdx = $targetIdx - 1\n\nif ($priorIdx -ge 1) {\n    # 1) Split a new (empty) paragraph mark in right after the prior\n    #    paragraph. This new mark inherits hint=\"eastAsia\".\n    $priorPara = $d.Paragraphs($priorIdx)\n    $priorPara.Range.InsertParagraphAfter()\n\n    # After the split, the structure is:\n    #   ... , priorIdx = prior paragraph,\n    #   priorIdx+1     = brand-new EMPTY paragraph (hint=\"eastAsia\"),\n    #   priorIdx+2     = original target paragraph (\"\u5f88\u82e6\", hint=\"default\",\n    #                     still carrying the _GoBack bookmark)\n    $newEastAsiaIdx = $priorIdx + 1\n    $oldTargetIdx = $priorIdx + 2\n\n    # 2) Move the target text into the freshly created eastAsia-hinted\n    #    paragraph.\n    $newEastAsiaPara = $d.Paragraphs($newEastAsiaIdx)\n    $newEastAsiaPara.Range.InsertAfter($targetText)\n\n    # 3) Clear the old paragraph's text (its mark \u2014 and the _GoBack\n    #    bookmark riding on it \u2014 stays put, still hint=\"default\") and\n    #    retype it with the new sentence. Because this mark is untouched\n    #    (not newly created), its hint stays \"default\", exactly matching\n    #    the freshly-authored paragraph in the target document.\n    $oldTargetPara = $d.Paragraphs($oldTargetIdx)\n    $oldTargetPara.Range.Text = \"\"\n    $oldTargetPara = $d.Paragraphs($oldTargetIdx)\n    $oldTargetPara.Range.InsertAfter($newText)\n} else {\n    # Target paragraph is the very first paragraph \u2014 no predecessor to\n    # borrow an \"eastAsia\" mark from; fall back to a plain append.\n    $targetPara = $d.Paragraphs($targetIdx)\n    $targetPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs($targetIdx + 1)\n    $newPara.Range.InsertAfter($newText)\n}\n\n# 4) Relocate the \"_GoBack\" bookmark to the very end of the document (the\n#    end of the freshly authored paragraph).\n$hadGoBack = $false\ntry {\n    $hadGoBack = $d.Bookmarks.Exists(\"_GoBack\")\n} catch {\n    $hadGoBack = $false\n}\n\n$docEnd = $d.Content\n$docEnd.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $docEnd) | Out-Null\n"}
